# "Implemented mthcutoff 6 months"
# - Adds a new "PSDNEW" candidate row (PSD ministry, posting start 1/1/2019)
#   to the "can (2)" sheet.
# - Updates the active sheet/selection/zoom bookkeeping to match where the
#   author was working when they made the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "pan (2)"
$ws2 = $wb.Worksheets.Item(2)   # "can (2)"

# --- Data edit: append row 17 to "can (2)", copying the format of row 16
#     (same pattern used by the existing "AAA"/"BBB" rows 15/16) then
#     overwriting the values.
$ws2.Range("A16").Copy()
$ws2.Range("A17").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("B16").Copy()
$ws2.Range("B17").PasteSpecial(-4122)

$ws2.Range("C16").Copy()
$ws2.Range("C17").PasteSpecial(-4122)

$ws2.Range("A17").Value = "PSDNEW"
$ws2.Range("B17").Value = "PSD"
$ws2.Range("C17").Value = 43466   # 1/1/2019

# --- View/selection bookkeeping ---

# Sheet1 ("pan (2)"): keep scrolled to row 22, zoom to 145%, move the
# selection to B18.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B18").Select()
$excel.ActiveWindow.Zoom = 145

# Sheet2 ("can (2)"): zoom to 130%, select C18, and leave it as the active
# (visible) sheet/tab.
$ws2.Activate()
$ws2.Range("C18").Select()
$excel.ActiveWindow.Zoom = 130

Write-Host "Applied mthcutoff 6 months edit"
